$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (GitHub Actions scheduled update).
# Columns: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
#
# Column D ("Price") holds plain text such as "1.00", "0.0000265" or
# "86.431.44" (dotted/grouped numbers). Excel would otherwise coerce these
# numeric-looking strings into real numbers (dropping trailing zeros,
# switching to scientific notation, etc.), so each Price cell we touch is
# explicitly formatted as Text ("@") before the value is written, keeping
# it identical in spirit to the original text content.

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value2 = '87.326.45'
$ws.Cells.Item(2,5).Value2 = '  +6.47%  '

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value2 = '3.294.19'
$ws.Cells.Item(3,5).Value2 = '  +3.34%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = '1.00'
$ws.Cells.Item(4,5).Value2 = '  +0.17%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = '213.72'
$ws.Cells.Item(5,5).Value2 = '  -0.84%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = '626.50'
$ws.Cells.Item(6,5).Value2 = '  +0.19%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = '0.388'
$ws.Cells.Item(7,5).Value2 = '  +34.70%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value2 = '0.999'
$ws.Cells.Item(8,5).Value2 = '  -0.01%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = '0.639'
$ws.Cells.Item(9,5).Value2 = '  +8.86%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = '3.293.23'
$ws.Cells.Item(10,5).Value2 = '  +3.39%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = '0.586'
$ws.Cells.Item(11,5).Value2 = '  -0.91%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = '0.0000265'
$ws.Cells.Item(12,5).Value2 = '  +2.28%  '

$ws.Cells.Item(13,5).Value2 = '  +5.58%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = '34.58'
$ws.Cells.Item(14,5).Value2 = '  +8.92%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = '3.898.16'
$ws.Cells.Item(15,5).Value2 = '  +3.28%  '

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value2 = '5.26'
$ws.Cells.Item(16,5).Value2 = '  -1.17%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value2 = '86.834.18'
$ws.Cells.Item(17,5).Value2 = '  +6.47%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value2 = '3.277.73'
$ws.Cells.Item(18,5).Value2 = '  +2.68%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value2 = '14.26'
$ws.Cells.Item(19,5).Value2 = '  +1.53%  '

# Rows 20/21 swap rank order: SuiNetwork moves above Uniswap.
$ws.Cells.Item(20,2).Value2 = 'SuiNetwork'
$ws.Cells.Item(20,3).Value2 = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = '2.98'
$ws.Cells.Item(20,5).Value2 = '  -7.13%  '

$ws.Cells.Item(21,2).Value2 = 'Uniswap'
$ws.Cells.Item(21,3).Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = '9.22'
$ws.Cells.Item(21,5).Value2 = '  +2.84%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = '438.62'
$ws.Cells.Item(22,5).Value2 = '  +0.61%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = '5.38'
$ws.Cells.Item(23,5).Value2 = '  +4.67%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = '7.26'
$ws.Cells.Item(24,5).Value2 = '  +0.08%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = '5.23'
$ws.Cells.Item(25,5).Value2 = '  -2.30%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value2 = '12.20'
$ws.Cells.Item(26,5).Value2 = '  +11.04%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = '3.473.32'
$ws.Cells.Item(27,5).Value2 = '  +3.65%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value2 = '77.01'
$ws.Cells.Item(28,5).Value2 = '  +0.34%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = '0.0000131'
$ws.Cells.Item(29,5).Value2 = '  +5.99%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = '1.00'
$ws.Cells.Item(30,5).Value2 = '  +0.27%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value2 = '0.176'
$ws.Cells.Item(31,5).Value2 = '  +26.97%  '

$ws.Cells.Item(32,5).Value2 = '  +0.75%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = '8.94'
$ws.Cells.Item(33,5).Value2 = '  -1.53%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = '555.96'
$ws.Cells.Item(34,5).Value2 = '  -5.05%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = '1.46'
$ws.Cells.Item(35,5).Value2 = '  -3.83%  '

$ws.Cells.Item(36,5).Value2 = '  -0.98%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = '6.91'
$ws.Cells.Item(37,5).Value2 = '  +12.25%  '

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value2 = '0.139'
$ws.Cells.Item(38,5).Value2 = '  -10.70%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = '22.72'
$ws.Cells.Item(39,5).Value2 = '  -0.45%  '

# Rows 40/41 swap rank order: FirstDigitalUSD moves above WhiteBITCoin.
$ws.Cells.Item(40,2).Value2 = 'FirstDigitalUSD'
$ws.Cells.Item(40,3).Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = '1.00'
$ws.Cells.Item(40,5).Value2 = '  +0.34%  '

$ws.Cells.Item(41,2).Value2 = 'WhiteBITCoin'
$ws.Cells.Item(41,3).Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = '21.74'
$ws.Cells.Item(41,5).Value2 = '  +4.53%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = '0.403'
$ws.Cells.Item(42,5).Value2 = '  -1.54%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value2 = '2.02'
$ws.Cells.Item(43,5).Value2 = '  -1.60%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value2 = '2.98'
$ws.Cells.Item(44,5).Value2 = '  -3.30%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = '1.00'
$ws.Cells.Item(45,5).Value2 = '  +0.04%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = '153.49'
$ws.Cells.Item(46,5).Value2 = '  -4.66%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = '181.28'
$ws.Cells.Item(47,5).Value2 = '  -3.59%  '

# Rows 48/49 swap rank order: OKB moves above ImmutableX.
$ws.Cells.Item(48,2).Value2 = 'OKB'
$ws.Cells.Item(48,3).Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value2 = '45.24'
$ws.Cells.Item(48,5).Value2 = '  +1.17%  '

$ws.Cells.Item(49,2).Value2 = 'ImmutableX'
$ws.Cells.Item(49,3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = '1.36'
$ws.Cells.Item(49,5).Value2 = '  +1.29%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = '4.26'
$ws.Cells.Item(50,5).Value2 = '  +1.26%  '

# Row 51: ARBITRUM drops off the list, replaced by Mantle.
$ws.Cells.Item(51,2).Value2 = 'Mantle'
$ws.Cells.Item(51,3).Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = '0.749'
$ws.Cells.Item(51,5).Value2 = '  -3.34%  '
